# Update column G ("K") values for rows 2-43 on Sheet1.
# These values reflect switching the source data generation from the
# old "Strike#" based computation to the new "K" computation, with
# regenerated std/mean and s_vals feeding into the computed K counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("G11").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("G14").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("G16").Value = 3
$ws.Range("G17").Value = 1
$ws.Range("G18").Value = 3
$ws.Range("G19").Value = 2
$ws.Range("G20").Value = 1
$ws.Range("G21").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("G23").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("G27").Value = 1
$ws.Range("G28").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("G31").Value = 1
$ws.Range("G32").Value = 3
$ws.Range("G33").Value = 2
$ws.Range("G34").Value = 0
$ws.Range("G35").Value = 3
$ws.Range("G37").Value = 2
$ws.Range("G39").Value = 1
$ws.Range("G40").Value = 1
$ws.Range("G41").Value = 0
$ws.Range("G43").Value = 1
